$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on Price (D) and Volume (E) columns so that
# numeric-looking strings (e.g. "30.417.79", "0.4740") are preserved exactly as text
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.417.79"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.868.34"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "245.92"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "0.4740"
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").Value = "0.2900"
$ws.Range("E8").Value = "  +1.64%  "
$ws.Range("D9").Value = "0.06491"
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("D10").Value = "21.94"
$ws.Range("E10").Value = "  +5.83%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "0.07711"
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("B12").Value = "Litecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D12").Value = "97.65"
$ws.Range("E12").Value = "  +3.87%  "
$ws.Range("D13").Value = "0.7350"
$ws.Range("E13").Value = "  +7.93%  "
$ws.Range("D14").Value = "1.870.62"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").Value = "5.116"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").Value = "273.60"
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("D17").Value = "30.403.59"
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("D19").Value = "0.000007553"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D21").Value = "2.114.83"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "5.224"
$ws.Range("E23").Value = "  +1.00%  "
$ws.Range("D24").Value = "6.157"
$ws.Range("E24").Value = "  +1.02%  "
$ws.Range("D25").Value = "9.267"
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("D26").Value = "164.07"
$ws.Range("E26").Value = "  -0.92%  "
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("D28").Value = "1.926"
$ws.Range("E28").Value = "  +2.19%  "
$ws.Range("D29").Value = "0.1000"
$ws.Range("E29").Value = "  +1.72%  "
$ws.Range("D30").Value = "1.366"
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("D31").Value = "1.507"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").Value = "4.303"
$ws.Range("E32").Value = "  +1.84%  "
$ws.Range("D33").Value = "4.156"
$ws.Range("E33").Value = "  +4.59%  "
$ws.Range("D34").Value = "0.04828"
$ws.Range("E34").Value = "  +2.62%  "
$ws.Range("D35").Value = "1.118"
$ws.Range("E35").Value = "  +0.74%  "
$ws.Range("D36").Value = "0.6960"
$ws.Range("E36").Value = "  +1.55%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "1.000"
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").Value = "2.713"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01855"
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.747"
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "6.295"
$ws.Range("E41").Value = "  -1.37%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "72.70"
$ws.Range("E42").Value = "  +3.36%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "1.965"
$ws.Range("E43").Value = "  +4.37%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.4184"
$ws.Range("E44").Value = "  +3.07%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "0.8338"
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "102.46"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "9.209"
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "7.009"
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "35.30"
$ws.Range("E50").Value = "  +2.51%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "926.69"
$ws.Range("E51").Value = "  +0.16%  "
